$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics (column C) and p-values (column D)
# for rows 2-11, per "Correcion a Diebold Mariano y revision de Cap1"

$updates = @(
    @{ Row = 2;  C = -0.7093769247427983; D = 0.4855414028498246 },
    @{ Row = 3;  C = -0.5714460091684388; D = 0.5734897604977247 },
    @{ Row = 4;  C = 1.668582464715588;   D = 0.1093715646692899 },
    @{ Row = 5;  C = 0.5390548527508098;  D = 0.5952626418708618 },
    @{ Row = 6;  C = 0.2211409009189114;  D = 0.8270223212614127 },
    @{ Row = 7;  C = 2.212980626857715;   D = 0.03757897336969873 },
    @{ Row = 8;  C = 1.039332685913267;   D = 0.309938551999825 },
    @{ Row = 9;  C = 2.562143847588255;   D = 0.01777144486479121 },
    @{ Row = 10; C = 1.4586859003467;     D = 0.1587820250708607 },
    @{ Row = 11; C = -1.38577060929277;   D = 0.1797053147528764 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
